$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Auto-fill the COM Port combo box value with "COM3" (text, not numeric)
$ws.Range("E2").Value = "COM3"

# Update the active selection to E2 (reflecting the cell just filled)
$ws.Range("E2").Select()
